$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rubrica")

$ws.Range("D6").Value = "PENDIENTE PEGAR EL WORD"
$ws.Range("C17").Value = "si"
$ws.Range("C20").Value = "si"
$ws.Range("C21").Value = "si"
$ws.Range("D7").Select()
